# Apply the "Add files via upload" edit:
#  - Trim trailing space/dot off a few header labels on Sheet1 (shared-string edits)
#  - Re-point Sheet1's view (zoom + selection)
#  - Populate the previously-empty Sheet2 with a header row + 3 rows of sample data,
#    including a new "right + vertical-center" cell style for the first data row
#  - Re-point Sheet2's view (selection)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------------
# Sheet1: clean up the three mistyped labels (trailing space / trailing dot).
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "Sr_No"
$ws1.Range("A3").Value = "Name"
$ws1.Range("A4").Value = "Work"
Write-Output "Sheet1 labels cleaned up"

# Sheet1 view: lower the crazy 400% zoom and move the active selection.
$ws1.Activate()
$excel.ActiveWindow.Zoom = 325
$ws1.Range("D3").Select()
Write-Output "Sheet1 view updated"

# ---------------------------------------------------------------------------
# Sheet2: fill in the header row + 3 rows of sample data.
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "Sr_No"
$ws2.Range("B1").Value = "Name"
$ws2.Range("C1").Value = "Work"
$ws2.Range("D1").Value = "Education"
# Header uses the same center/center alignment already used on Sheet1's header.
$ws2.Range("A1:D1").HorizontalAlignment = -4108   # xlCenter
$ws2.Range("A1:D1").VerticalAlignment = -4108     # xlCenter
Write-Output "Sheet2 header written"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 3
$ws2.Range("D2").Value = 4

# Build the new right/vertical-center style cleanly on a scratch cell first
# (elsewhere, on Sheet3) and copy *just the format* across — doing the two
# alignment assignments directly on a multi-cell range leaves a stray unused
# cellXf behind, so routing through a single-cell seed + format-only paste
# keeps the style table minimal, matching a single freshly added xf.
$seed = $ws3.Range("Z100")
$seed.HorizontalAlignment = -4152   # xlRight
$seed.VerticalAlignment = -4108     # xlCenter
$seed.Copy()
$ws2.Range("A2:D2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$seed.Clear()
Write-Output "Sheet2 row 2 styled"

$ws2.Range("A3").Value = 5
$ws2.Range("B3").Value = 6
$ws2.Range("C3").Value = 7
$ws2.Range("D3").Value = 8

$ws2.Range("A4").Value = 9
$ws2.Range("B4").Value = 10
$ws2.Range("C4").Value = 11
$ws2.Range("D4").Value = 12
Write-Output "Sheet2 data rows written"

# Sheet2 view: move the active selection.
$ws2.Activate()
$ws2.Range("H8").Select()
Write-Output "Sheet2 view updated"

# Leave Sheet1 as the active/selected tab, as in the target workbook.
$ws1.Activate()
Write-Output "done"
